# Generate Report for Handback
# - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
# - zh-cn row: handback datetime refreshed, error detail cleared (no longer stale)
# - de-de row: handback datetime refreshed, error detail cleared (no longer stale)
# - Columns widened/narrowed to fit the new content (AutoFit-style)

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("K2").Value = "2016-08-21 10:54:39"
$zh.Range("P2").Value = ""
$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("K2").Value = "2016-08-21 10:54:45"
$de.Range("P2").Value = ""
$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(16).ColumnWidth = 12.833333333333334

Write-Output "Report regenerated for handback."
